{"js": "// Office.js (Word JavaScript API) script.\n// Body of: async (context) => { ... }\n//\n// The \"Log\" table has a row documenting the quitGameCallback change\n// (\"Hi\u1ec7n th\u1ef1c h\u00e0m quitGameCallback trong class GameLayer\" / \"S\u01a1n\") that is\n// immediately followed by a still-empty row. This change fills that empty\n// row with a new log entry:\n//   6/10/2012 | B\u1ed5 sung h\u00e0m menuExitCallback cho class MenuScene | S\u01a1n\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\ntable.load(\"values\");\nawait context.sync();\n\n// Locate the \"quitGameCallback\" row dynamically, then target the row right\n// after it (the first still-empty log row) instead of hard-coding indexes.\nlet quitRowIndex = -1;\nfor (let i = 0; i < table.values.length; i++) {\n  const middleCol = table.values[i][1] || \"\";\n  if (middleCol.indexOf(\"quitGameCallback\") !== -1) {\n    quitRowIndex = i;\n    break;\n  }\n}\nif (quitRowIndex === -1) {\n  throw new Error(\"Could not find the quitGameCallback log row.\");\n}\nconst newRowIndex = quitRowIndex + 1;\n\n// The \"S\u01a1n\" author cell on the quitGameCallback row owns Word's internal\n// \"_GoBack\" bookmark (last-edit marker). Re-type it cleanly so the bookmark\n// moves with the edit, same as Word would do.\nconst authorCellAbove = table.getCell(quitRowIndex, 2);\nauthorCellAbove.body.clear();\nawait context.sync();\nauthorCellAbove.body.insertText(\"S\u01a1n\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Fill in the new log entry row.\nconst dateCell = table.getCell(newRowIndex, 0);\ndateCell.body.insertText(\"6/10/2012\", Word.InsertLocation.replace);\n\nconst descCell = table.getCell(newRowIndex, 1);\ndescCell.body.insertText(\n  \"B\u1ed5 sung h\u00e0m menuExitCallback cho class MenuScene\",\n  Word.InsertLocation.replace\n);\nawait context.sync();\n\nconst authorCell = table.getCell(newRowIndex, 2);\nauthorCell.body.insertText(\"S\u01a1n\", Word.InsertLocation.replace);\nawait context.sync();\n\n// Drop the \"_GoBack\" bookmark where Word would leave it: right after typing\n// \"menuE\" (mid-word), before continuing with \"xitCallback...\".\nconst searchResults = descCell.body.search(\"menuE\", { matchCase: true });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length > 0) {\n  const afterTyped = searchResults.items[0].getRange(Word.RangeLocation.after);\n  afterTyped.insertBookmark(\"_GoBack\");\n  await context.sync();\n}\n", "ps1": "# Word COM interop (PowerShell-style) script.\n# $word.ActiveDocument is already open as $d below.\n#\n# The \"Log\" table has a row documenting the quitGameCallback change\n# (\"Hien thuc ham quitGameCallback trong class GameLayer\" / \"Son\") that is\n# immediately followed by a still-empty row. This change fills that empty\n# row with a new log entry:\n#   6/10/2012 | Bo sung ham menuExitCallback cho class MenuScene | Son\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables(1)\n\n# Locate the \"quitGameCallback\" row dynamically, then target the row right\n# after it (the first still-empty log row) instead of hard-coding indexes.\n$quitRow = 0\nfor ($r = 1; $r -le $tbl.Rows.Count; $r++) {\n    $desc = $tbl.Cell($r, 2).Range.Text\n    if ($desc -like \"*quitGameCallback*\") {\n        $quitRow = $r\n        break\n    }\n}\nif ($quitRow -eq 0) {\n    throw \"Could not find the quitGameCallback log row.\"\n}\n$newRow = $quitRow + 1\n\n# The \"S\u01a1n\" author cell on the quitGameCallback row owns Word's internal\n# \"_GoBack\" bookmark (last-edit marker). Deleting the paragraph's range and\n# retyping it drops that bookmark, then we re-insert the text cleanly -\n# mirroring how the bookmark moves with a fresh edit in Word.\n$authorCellAbove = $tbl.Cell($quitRow, 3)\n$paraAbove = $authorCellAbove.Range.Paragraphs(1)\n$paraAbove.Range.Delete()\n$authorCellAbove.Range.InsertAfter(\"S\u01a1n\")\n\n# Fill in the new log entry row.\n$tbl.Cell($newRow, 1).Range.Text = \"6/10/2012\"\n$descCell = $tbl.Cell($newRow, 2)\n$descCell.Range.Text = \"B\u1ed5 sung h\u00e0m menuExitCallback cho class MenuScene\"\n$tbl.Cell($newRow, 3).Range.Text = \"S\u01a1n\"\n\n# Drop the \"_GoBack\" bookmark where Word would leave it: right after typing\n# \"menuE\" (mid-word), before continuing with \"xitCallback...\".\n$findRange = $d.Range($descCell.Range.Start, $descCell.Range.End)\n$found = $findRange.Find.Execute(\"menuE\")\nif ($found) {\n    $bmRange = $d.Range($findRange.End, $findRange.End)\n    $d.Bookmarks.Add(\"_GoBack\", $bmRange)\n}\n"}
